$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update the URL (pythia -> cicada)
$ws.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/vaccine-type"

# 2. Update the Date
$ws.Range("B8").Value = "2026-02-11T14:37:07-05:00"

# 3. Insert a new "Jurisdiction" row after "Contact" (row 10), before "Description" (row 11)
#    Copy the format from the row that will end up below it so the new row matches the
#    existing data-row style instead of Excel's blank default style.
$ws.Rows.Item(11).Copy()
$ws.Rows.Item(11).Insert()
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
